$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.480.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.91%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.117.55"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.34%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +13.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.69"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.77"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.85%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.647.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.92%  "

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.132.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.66%  "

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.456.60"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.36"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.25"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.08"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.41%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.16%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.85"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.06"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0866"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.36%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +14.05%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "448.33"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.86%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0371"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.894.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.280"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.52%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.52"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.74"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.67"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.30%  "
